$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = [double]"0.000209991721317872"
$ws.Range("E2").Value = [double]"0.000209991721317872"

# Row 3
$ws.Range("D3").Value = [double]"0.9998905564803414"
$ws.Range("E3").Value = [double]"0.9998905564803414"

# Row 4
$ws.Range("D4").Value = [double]"0.00105172623047923"
$ws.Range("E4").Value = [double]"0.00105172623047923"

# Row 5
$ws.Range("D5").Value = [double]"1.040333716405986E-08"
$ws.Range("E5").Value = [double]"1.040333716405986E-08"

# Row 6
$ws.Range("D6").Value = [double]"0.05081521290413182"
$ws.Range("E6").Value = [double]"0.05081521290413182"

# Row 7
$ws.Range("D7").Value = [double]"0.999999999942754"
$ws.Range("E7").Value = [double]"5.72459857295371E-11"

# Row 8
$ws.Range("C8").Value = $false
$ws.Range("D8").Value = [double]"7.302302926020122E-11"
$ws.Range("E8").Value = [double]"0.999999999926977"

# Row 9
$ws.Range("C9").Value = $false
$ws.Range("D9").Value = [double]"0.01007019350147272"
$ws.Range("E9").Value = [double]"0.9899298064985272"

# Row 10
$ws.Range("C10").Value = $false
$ws.Range("D10").Value = [double]"1.627354027690754E-08"
$ws.Range("E10").Value = [double]"0.9999999837264597"

# Row 11
$ws.Range("D11").Value = [double]"0.9999981542118342"
$ws.Range("E11").Value = [double]"1.845788165755202E-06"
$ws.Range("F11").Value = [double]"5.504566669464111"
$ws.Range("G11").Value = [double]"0.6"

# Row 12
$ws.Range("D12").Value = [double]"6.907648414939411E-07"
$ws.Range("E12").Value = [double]"6.907648414939411E-07"

# Row 13
$ws.Range("D13").Value = [double]"0.9999999997789986"
$ws.Range("E13").Value = [double]"0.9999999997789986"

# Row 14
$ws.Range("D14").Value = [double]"8.075888261925175E-05"
$ws.Range("E14").Value = [double]"8.075888261925175E-05"

# Row 15
$ws.Range("D15").Value = [double]"6.829117625862511E-09"
$ws.Range("E15").Value = [double]"6.829117625862511E-09"

# Row 16
$ws.Range("D16").Value = [double]"0.03304493395354943"
$ws.Range("E16").Value = [double]"0.03304493395354943"

# Row 17
$ws.Range("D17").Value = [double]"1"
$ws.Range("E17").Value = [double]"0"

# Row 18
$ws.Range("C18").Value = $false
$ws.Range("D18").Value = [double]"7.444714747088096E-17"
$ws.Range("E18").Value = [double]"0.9999999999999999"

# Row 19
$ws.Range("C19").Value = $false
$ws.Range("D19").Value = [double]"0.0001300040277642591"
$ws.Range("E19").Value = [double]"0.9998699959722357"

# Row 20
$ws.Range("C20").Value = $false
$ws.Range("D20").Value = [double]"1.528929464234998E-12"
$ws.Range("E20").Value = [double]"0.9999999999984711"

# Row 21
$ws.Range("D21").Value = [double]"0.9999999609086931"
$ws.Range("E21").Value = [double]"3.909130685819662E-08"
$ws.Range("F21").Value = [double]"9.55573844909668"
$ws.Range("G21").Value = [double]"0.6"
